# banned_flight_599H.xlsx edit:
#  - "Singapore Airlines" -> "SingaporeAirlines" (airline_e column)
#  - "Singapore " (trailing space) -> "Singapore" (land_e / city_e columns)
#  - sheet view scrolled/selection moved (topLeftCell A31->A10, selection J62->E44)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix "Singapore Airlines" -> "SingaporeAirlines" in column B (airline_e) ---
$airlineCells = @("B23", "B44", "B45", "B56")
foreach ($addr in $airlineCells) {
    $cell = $ws.Range($addr)
    if ($cell.Value2 -eq "Singapore Airlines") {
        $cell.Value = "SingaporeAirlines"
    }
}

# --- Fix "Singapore " (trailing space) -> "Singapore" in land_e / city_e columns ---
$singaporeCells = @("E44", "G44", "E45", "G45", "E49", "G49", "E56", "G56")
foreach ($addr in $singaporeCells) {
    $cell = $ws.Range($addr)
    if ($cell.Value2 -eq "Singapore ") {
        $cell.Value = "Singapore"
    }
}

# --- Update the active sheet view's scroll position / selection ---
$ws.Range("A10").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E44").Select()
